# Apply the BreakoutBoard BOM edits via Excel COM interop
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the unused empty sheets (Sheet2, Sheet3)
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()
$excel.DisplayAlerts = $true

# 2. Rename the remaining sheet
$ws.Name = "Breakout Board - Rev B"

# 3. Update the title text in A1 to the new BOM title
$ws.Range("A1").Value = "Bill of Materials for 'Marmote - Breakout Board Rev B (Smoky)'"

# 4. Swap the "Item #" numbering between row 4 (J1) and row 5 (J2)
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 2

# 5. Update the selection to the header row range, no distinct active cell
$ws.Range("A1:L1").Select()
